# Rename header labels on the existing sheets
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# Reuse the existing header style (bold, bordered, centered) from sheet 1
$ws1.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Reuse the existing date-formatted style from sheet 1 for column A
$ws1.Range("A2").Copy()
$newSheet.Range("A2:A31").PasteSpecial(-4122)

# Header row
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Forecast data rows
$data = @(
    @(45424.99999999999, 10, -196.4145677260732, 220.8566418841556),
    @(45445.99999999999, 39, -163.9269254159585, 241.9659040061348),
    @(45459.99999999999, 59, -150.2025348828862, 253.2027074527664),
    @(45466.99999999999, 69, -136.3968505570568, 275.7461264762698),
    @(45473.99999999999, 79, -112.9611350644414, 279.8849654522818),
    @(45480.99999999999, 89, -109.1933111885373, 290.215545696084),
    @(45487.99999999999, 99, -101.9208660981232, 295.9917215672568),
    @(45494.99999999999, 109, -88.36641911038545, 308.0630981958453),
    @(45501.99999999999, 119, -82.39046591259702, 319.5999590472888),
    @(45508.99999999999, 129, -78.25884518710109, 331.8051205372478),
    @(45515.99999999999, 139, -44.67327122042853, 342.4270824043487),
    @(45522.99999999999, 148, -43.57511306051371, 358.5619984962117),
    @(45536.99999999999, 168, -45.9299993071162, 374.4981646999047),
    @(45543.99999999999, 178, -10.55052910438392, 391.2151477716372),
    @(45550.99999999999, 188, -8.10037613946108, 391.2386694095987),
    @(45557.99999999999, 198, -14.51674526567687, 403.0252945390952),
    @(45564.99999999999, 208, 2.436879178163244, 428.1418666637862),
    @(45571.99999999999, 218, 17.0741825892231, 434.4577427677601),
    @(45578.99999999999, 228, 29.16798457483699, 435.2370977702468),
    @(45585.99999999999, 238, 40.3598965073331, 444.7247999183609),
    @(45592.99999999999, 248, 58.36235793536802, 459.8290903116522),
    @(45634.99999999999, 307, 90.87784540427924, 506.3168702415031),
    @(45641.99999999999, 317, 111.1018713305163, 511.6930183611883),
    @(45648.99999999999, 327, 124.1680371878644, 543.7510985447611),
    @(45655.99999999999, 337, 129.1493006766678, 536.910939022645),
    @(45662.99999999999, 347, 148.2566680194989, 536.3937688406257),
    @(45669.99999999999, 357, 158.8647418074667, 549.0643895445734),
    @(45676.99999999999, 367, 160.8321644225105, 572.3739253567779),
    @(45683.99999999999, 377, 169.7638107067497, 575.5674606633975),
    @(45690.99999999999, 386, 172.9288405317371, 593.7043788543224)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $r++
}
